$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Training Data Issue (#48):
# Column BF ("Date") held a malformed date string ("6-18-2011-12") on
# every data row because of the way the NBA stats source rendered the
# date - it needs to read "2012-06-18" instead (the data was off by a
# day / wrongly formatted). Correct BF2:BF31 in place.
#
# Plain assignment of an ISO-looking string to .Value auto-converts the
# cell to a real date serial, which is not what the source data has
# (it stays a literal text string). Stamping the cell as Text first
# keeps the assignment literal, then the style is restored to Normal so
# the cell's formatting is left exactly as it was.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    $current = $cell.Value2()
    if ($current -eq "6-18-2011-12") {
        $cell.NumberFormat = "@"
        $cell.Value = "2012-06-18"
        $cell.Style = "Normal"
    }
}
